# Applies the edit described by the diff: several rows in the "Artfynd"
# sheet had their per-observation fields (Id, Ost, Nord, Starttid, Sluttid
# - and in two cases the whole row) shuffled between rows.
#
# The affected row groups (1-based worksheet rows, full column span A:AY):
#   - Rows 5,6,7,8   : rotate so new5<-old6, new6<-old7, new7<-old8, new8<-old5
#   - Rows 20,21     : swap entire rows
#   - Rows 25,26     : swap entire rows
#   - Rows 37,38,39  : rotate so new37<-old39, new38<-old37, new39<-old38

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstCol = "A"
$lastCol  = "AY"

function Get-RowValues($rowNum) {
    return $ws.Range("$firstCol$rowNum`:$lastCol$rowNum").Value2
}

function Set-RowValues($rowNum, $vals) {
    $ws.Range("$firstCol$rowNum`:$lastCol$rowNum").Value2 = $vals
}

# --- Rotate rows 5,6,7,8 (new[r] = old[r+1], wrapping 8 -> 5) ---
$r5 = Get-RowValues 5
$r6 = Get-RowValues 6
$r7 = Get-RowValues 7
$r8 = Get-RowValues 8

Set-RowValues 5 $r6
Set-RowValues 6 $r7
Set-RowValues 7 $r8
Set-RowValues 8 $r5

# --- Swap rows 20 and 21 ---
$r20 = Get-RowValues 20
$r21 = Get-RowValues 21

Set-RowValues 20 $r21
Set-RowValues 21 $r20

# --- Swap rows 25 and 26 ---
$r25 = Get-RowValues 25
$r26 = Get-RowValues 26

Set-RowValues 25 $r26
Set-RowValues 26 $r25

# --- Rotate rows 37,38,39 (new37<-old39, new38<-old37, new39<-old38) ---
$r37 = Get-RowValues 37
$r38 = Get-RowValues 38
$r39 = Get-RowValues 39

Set-RowValues 37 $r39
Set-RowValues 38 $r37
Set-RowValues 39 $r38
